# Generate Report for Handoff
#
# Updates the localization-status report after a fresh handoff-xliff
# generation run:
#   - Priority for the three still-pending files (7fc7faea, c3732db0,
#     c647cd9d, feb26b48) flips from "low" to "ht" on both the zh-cn and
#     de-de language sheets.
#   - The zh-cn "Latest Handoff Datetime" for those same rows advances to
#     the new generation timestamp.
#   - The shared "Latest HO Xliff Generate Date" (Overview sheet, which is
#     also shown as the de-de "Latest Handoff Datetime") advances too.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 4-7
$overview.Range("G4:G7").Value = "2016-08-27 20:31:30"

# zh-cn sheet: Priority (E) and Latest Handoff Datetime (H) for rows 4-7
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-27 20:31:24"

# de-de sheet: Priority (E) for rows 4-7
$dede.Range("E4:E7").Value = "ht"

# de-de sheet "Latest Handoff Datetime" (H) mirrors the same refreshed
# generation timestamp as the Overview sheet's G4:G7 above.
$dede.Range("H4:H7").Value = "2016-08-27 20:31:30"
